$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1

$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.Item(8).LineStyle = 1
$ws1.Range("D1").Borders.Item(9).LineStyle = 1
$ws1.Range("D1").Borders.Item(10).LineStyle = 1

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2.Range("C1").ClearFormats()
$ws2.Range("C1").Borders.Item(8).LineStyle = 1
$ws2.Range("C1").Borders.Item(9).LineStyle = 1

$ws2.Range("D1").ClearFormats()
$ws2.Range("D1").Borders.Item(8).LineStyle = 1
$ws2.Range("D1").Borders.Item(9).LineStyle = 1
$ws2.Range("D1").Borders.Item(10).LineStyle = 1

$ws2.Range("F1").ClearFormats()
$ws2.Range("F1").Borders.Item(8).LineStyle = 1
$ws2.Range("F1").Borders.Item(9).LineStyle = 1

$ws2.Range("G1").ClearFormats()
$ws2.Range("G1").Borders.Item(8).LineStyle = 1
$ws2.Range("G1").Borders.Item(9).LineStyle = 1
$ws2.Range("G1").Borders.Item(10).LineStyle = 1

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
